$d = $word.ActiveDocument

# 1. "positions in the IT level from middle+" -> "positions in IT level from middle+"
$d.Content.Find.Execute(
    "positions in the IT level from middle+",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "positions in IT level from middle+", 2
)

# 2. "Zabbix, Jira administration" -> "Zabbix, Hujinn, Jira administration"
$d.Content.Find.Execute(
    "Zabbix, Jira administration",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Zabbix, Hujinn, Jira administration", 2
)

# 3. "I have the talent to install" -> "I have a talent to install"
$d.Content.Find.Execute(
    "I have the talent to install",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "I have a talent to install", 2
)
